# Update the "Organization website (if available)" cell (A10/B10 on
# "Лист 1") from the old domain to the new one, and leave the selection
# on B9 (matching the saved cursor position in the updated workbook).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист 1")

$ws.Range("B10").Value = "www.stat.gov.kg"

$ws.Activate()
$ws.Range("B9").Select()
